# Regenerate merged AHB files
# -----------------------------------------------------------------------
# This script:
#   1. Renames the 20 header labels in row 1 from "<name>_old" / "<name>_new"
#      to "<name>_FV2410" / "<name>_FV2504" (the "diff" header stays as is).
#   2. Turns the data range A1:U85 into an Excel Table ("Table1") without
#      picking up a table style / header dxf (so the rest of the workbook's
#      styles stay untouched), while keeping the header row's original
#      formatting (bold, gray fill, borders, centered + wrap text).
#   3. Freezes the header row (pane split after row 1).
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# Step 1: rename the header cells (A1:J1 = "_old" -> "_FV2410",
#         L1:U1 = "_new" -> "_FV2504"; K1 = "diff" is unchanged)
# -----------------------------------------------------------------------
$ws.Range("A1").Value2 = "Segmentname_FV2410"
$ws.Range("B1").Value2 = "Segmentgruppe_FV2410"
$ws.Range("C1").Value2 = "Segment_FV2410"
$ws.Range("D1").Value2 = "Datenelement_FV2410"
$ws.Range("E1").Value2 = "Segment ID_FV2410"
$ws.Range("F1").Value2 = "Code_FV2410"
$ws.Range("G1").Value2 = "Qualifier_FV2410"
$ws.Range("H1").Value2 = "Beschreibung_FV2410"
$ws.Range("I1").Value2 = "Bedingungsausdruck_FV2410"
$ws.Range("J1").Value2 = "Bedingung_FV2410"

$ws.Range("L1").Value2 = "Segmentname_FV2504"
$ws.Range("M1").Value2 = "Segmentgruppe_FV2504"
$ws.Range("N1").Value2 = "Segment_FV2504"
$ws.Range("O1").Value2 = "Datenelement_FV2504"
$ws.Range("P1").Value2 = "Segment ID_FV2504"
$ws.Range("Q1").Value2 = "Code_FV2504"
$ws.Range("R1").Value2 = "Qualifier_FV2504"
$ws.Range("S1").Value2 = "Beschreibung_FV2504"
$ws.Range("T1").Value2 = "Bedingungsausdruck_FV2504"
$ws.Range("U1").Value2 = "Bedingung_FV2504"

# -----------------------------------------------------------------------
# Step 2: convert A1:U85 into a native Excel table.
#
# Creating a ListObject on top of an already-formatted header row makes
# the engine capture the existing header formatting as a new dxf
# (headerRowDxfId) and register a default table style, which would touch
# xl/styles.xml even though the header cells end up looking the same.
# To avoid introducing that noise we temporarily park a copy of the
# header formatting in an unused scratch row, strip the formatting from
# the real header row, build the table (so it doesn't "see" any format
# to remember), then paste the original formatting back onto the header
# row and remove the scratch row again.
# -----------------------------------------------------------------------
$scratchRow = $ws.Rows.Item(1000)
$scratch = $ws.Range("A1000:U1000")

$headerRange = $ws.Range("A1:U1")
$headerRange.Copy() | Out-Null
$scratch.PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$headerRange.ClearFormats() | Out-Null

$rng = $ws.Range("A1:U85")
$lo = $ws.ListObjects.Add(1, $rng, [System.Reflection.Missing]::Value, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

$scratch.Copy() | Out-Null
$headerRange.PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = $false
$scratchRow.Delete() | Out-Null

# -----------------------------------------------------------------------
# Step 3: freeze the header row.
# -----------------------------------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

Write-Host "done"
